$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "max" (D) and "min" (E) values for row 2
$ws.Range("D2").Value = 34
$ws.Range("E2").Value = 56

# Fix the "max" (D) and "min" (E) values for row 9
$ws.Range("D9").Value = 32
$ws.Range("E9").Value = 45

# Update the active selection to match the edit
$ws.Range("E9").Select()
